$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# statut_name values ("B" column) are missing a space before the colon.
# Add the missing space, matching each row's "statut" code in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $statut = $ws.Cells.Item($r, 1).Value()
    if ($statut -eq "4") {
        $ws.Cells.Item($r, 2).Value = "4 : pas de résultats postés ni publiés"
    }
    elseif ($statut -eq "1") {
        $ws.Cells.Item($r, 2).Value = "1 : résultats postés ou publiés dans les 12 mois"
    }
}
